$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 208.42857
$ws.Range("I5").Value = 212.25
$ws.Range("J5").Value = 203.33333
$ws.Range("K5").Value = 212.25
$ws.Range("L5").Value = 203.33333
$ws.Range("M5").Value = -97.25
$ws.Range("N5").Value = -433.33333
$ws.Range("H17").Value = 3214.6738
$ws.Range("J17").Value = 3383.1396
$ws.Range("L17").Value = 10149.4188
$ws.Range("N17").Value = -10485.4188
$ws.Range("H70").Value = 5423.9434
$ws.Range("J70").Value = 5899.5957
$ws.Range("L70").Value = 17698.7871
$ws.Range("N70").Value = -18238.7871
$ws.Range("H73").Value = 5423.9434
$ws.Range("J73").Value = 5899.5957
$ws.Range("L73").Value = 17698.7871
$ws.Range("N73").Value = -19570.7871
$ws.Range("H106").Value = 2656.5
$ws.Range("I106").Value = 2807.0908
$ws.Range("K106").Value = 2807.0908
$ws.Range("M106").Value = -2176.0908
$ws.Range("H138").Value = 3348.2373
$ws.Range("J138").Value = 3376.1914
$ws.Range("L138").Value = 10128.5742
$ws.Range("N138").Value = -20408.5742

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 3760.75
$ws.Range("J15").Value = 3760.75
$ws.Range("L15").Value = 3760.75
$ws.Range("N15").Value = -4460.75
$ws.Range("H23").Value = 17500
$ws.Range("I23").Value = 17500
$ws.Range("K23").Value = 17500
$ws.Range("M23").Value = -17241
$ws.Range("H28").Value = 2287
$ws.Range("J28").Value = 1460
$ws.Range("L28").Value = 1460
$ws.Range("N28").Value = -1844
$ws.Range("H32").Value = 4325.25
$ws.Range("I32").Value = 2752.7715
$ws.Range("K32").Value = 2752.7715
$ws.Range("M32").Value = -2465.7715
$ws.Range("H45").Value = 3998316.5
$ws.Range("I45").Value = 5534212
$ws.Range("K45").Value = 5534212
$ws.Range("M45").Value = -5533835
$ws.Range("H61").Value = 3245
$ws.Range("J61").Value = 3492.5
$ws.Range("L61").Value = 3492.5
$ws.Range("N61").Value = -3916.5
$ws.Range("H99").Value = 2287
$ws.Range("J99").Value = 1460
$ws.Range("L99").Value = 1460
$ws.Range("N99").Value = -7450
$ws.Range("H132").Value = 1759.3617
$ws.Range("J132").Value = 2195.1333
$ws.Range("L132").Value = 6585.3999
$ws.Range("N132").Value = -11645.3999
$ws.Range("H136").Value = 3245
$ws.Range("J136").Value = 3492.5
$ws.Range("L136").Value = 10477.5
$ws.Range("N136").Value = -15577.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 1493.3334
$ws.Range("I37").Value = 2900
$ws.Range("J37").Value = 790
$ws.Range("K37").Value = 2900
$ws.Range("L37").Value = 790
$ws.Range("M37").Value = -2763
$ws.Range("N37").Value = -1064
$ws.Range("H86").Value = 4552930
$ws.Range("J86").Value = 3099.3333
$ws.Range("L86").Value = 3099.3333
$ws.Range("N86").Value = -5345.3333
$ws.Range("H89").Value = 4552930
$ws.Range("J89").Value = 3099.3333
$ws.Range("L89").Value = 15496.6665
$ws.Range("N89").Value = -26728.6665
$ws.Range("H99").Value = 4497543
$ws.Range("I99").Value = 5756091.5
$ws.Range("K99").Value = 5756091.5
$ws.Range("M99").Value = -5754593.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37158.42
$ws.Range("I31").Value = 1409.4615
$ws.Range("K31").Value = 1409.4615
$ws.Range("M31").Value = -1114.4615
$ws.Range("H34").Value = 37158.42
$ws.Range("I34").Value = 1409.4615
$ws.Range("K34").Value = 1409.4615
$ws.Range("M34").Value = -1207.4615
$ws.Range("H50").Value = 4157.8945
$ws.Range("J50").Value = 4157.8945
$ws.Range("L50").Value = 4157.8945
$ws.Range("N50").Value = -5407.8945
$ws.Range("H51").Value = 21049
$ws.Range("J51").Value = 32099
$ws.Range("L51").Value = 32099
$ws.Range("N51").Value = -33571
$ws.Range("H59").Value = 40000
$ws.Range("J59").Value = 40000
$ws.Range("L59").Value = 40000
$ws.Range("N59").Value = -42290
$ws.Range("H60").Value = 198.33333
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("H61").Value = 21049
$ws.Range("J61").Value = 32099
$ws.Range("L61").Value = 32099
$ws.Range("N61").Value = -32795
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("H97").Value = 24998.75
$ws.Range("J97").Value = 24998.75
$ws.Range("L97").Value = 24998.75
$ws.Range("N97").Value = -26980.75
$ws.Range("H134").Value = 3368
$ws.Range("I134").Value = 2570.7
$ws.Range("K134").Value = 7712.099999999999
$ws.Range("M134").Value = -5177.099999999999
$ws.Range("H140").Value = 67323
$ws.Range("J140").Value = 88984.5
$ws.Range("L140").Value = 88984.5
$ws.Range("N140").Value = -99344.5
$ws.Range("N60").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 67890.87
$ws.Range("I5").Value = 548.125
$ws.Range("K5").Value = 1644.375
$ws.Range("M5").Value = -1532.375
$ws.Range("H135").Value = 67890.87
$ws.Range("I135").Value = 548.125
$ws.Range("K135").Value = 4933.125
$ws.Range("M135").Value = -2398.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 35000
$ws.Range("J33").Value = 35000
$ws.Range("L33").Value = 35000
$ws.Range("N33").Value = -35504
$ws.Range("H122").Value = 386340.4
$ws.Range("I122").Value = 513944.34
$ws.Range("K122").Value = 1541833.02
$ws.Range("M122").Value = -1539383.02
$ws.Range("H132").Value = 3160.9033
$ws.Range("I132").Value = 2626.8333
$ws.Range("J132").Value = 4992
$ws.Range("K132").Value = 7880.499899999999
$ws.Range("L132").Value = 14976
$ws.Range("M132").Value = -5350.499899999999
$ws.Range("N132").Value = -20036

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 112345.75
$ws.Range("I22").Value = 178813.2
$ws.Range("J22").Value = 1566.6666
$ws.Range("K22").Value = 178813.2
$ws.Range("L22").Value = 1566.6666
$ws.Range("M22").Value = -178518.2
$ws.Range("N22").Value = -2156.6666
$ws.Range("H27").Value = 112345.75
$ws.Range("I27").Value = 178813.2
$ws.Range("J27").Value = 1566.6666
$ws.Range("K27").Value = 178813.2
$ws.Range("L27").Value = 1566.6666
$ws.Range("M27").Value = -178706.2
$ws.Range("N27").Value = -1780.6666
$ws.Range("H122").Value = 4455.433
$ws.Range("I122").Value = 2775.353
$ws.Range("K122").Value = 8326.059000000001
$ws.Range("M122").Value = -5876.059000000001
$ws.Range("H132").Value = 4640.04
$ws.Range("I132").Value = 3747.2942
$ws.Range("J132").Value = 6537.125
$ws.Range("K132").Value = 11241.8826
$ws.Range("L132").Value = 19611.375
$ws.Range("M132").Value = -8711.882599999999
$ws.Range("N132").Value = -24671.375
$ws.Range("H136").Value = 85032.08
$ws.Range("I136").Value = 121764.88
$ws.Range("J136").Value = 6974.875
$ws.Range("K136").Value = 365294.64
$ws.Range("L136").Value = 20924.625
$ws.Range("M136").Value = -362744.64
$ws.Range("N136").Value = -26024.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 35882.668
$ws.Range("J37").Value = 35882.668
$ws.Range("L37").Value = 35882.668
$ws.Range("N37").Value = -36288.668
$ws.Range("H107").Value = 33335468
$ws.Range("I107").Value = 45457184
$ws.Range("J107").Value = 746.25
$ws.Range("K107").Value = 136371552
$ws.Range("L107").Value = 2238.75
$ws.Range("M107").Value = -136369632
$ws.Range("N107").Value = -6078.75
$ws.Range("H122").Value = 1794.0416
$ws.Range("I122").Value = 1380.6666
$ws.Range("K122").Value = 4141.9998
$ws.Range("M122").Value = -1691.9998
$ws.Range("H136").Value = 2504.3333
$ws.Range("I136").Value = 1881.25
$ws.Range("K136").Value = 5643.75
$ws.Range("M136").Value = -3093.75
